# Add a second user row (id, username, password, name) to the users sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces these digit-only values to be stored as text,
# matching how the existing id column (e.g. 309044071) is stored.
$ws.Range("A2").Value = "'123456789"
$ws.Range("B2").Value = "raful9"
$ws.Range("C2").Value = "123456a!"
$ws.Range("D2").Value = "rafael"
